# "Generate Report for Archive"
#
# 1. The localization status for every tracked file flips from
#    "Ready for handoff" to "In Translation" (Overview!E/F and the
#    per-language sheets' Status column, all backed by the same shared
#    string, so a single find/replace across every sheet keeps every
#    occurrence in sync).
# 2. Because the new status text is shorter than the old one, the
#    Status-bearing columns (Overview!E:F, and column C on the zh-cn /
#    de-de sheets) got narrower when the report was regenerated/saved.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Target stored column width (the "width" attribute Excel writes into
# <col .../> in the sheet XML) after the narrower text was applied.
$targetStoredWidth = 13.4101845877511

# The ColumnWidth COM setter snaps to Excel's internal character grid
# (integer pixel widths built from the workbook's Maximum Digit Width,
# 6px here, plus 5px of cell padding: pixels = round(chars*MDW)+5, and
# the value actually stored back in the file is pixels/MDW). Solve that
# relationship in reverse so the ColumnWidth we feed in lands on the
# closest achievable grid point to the real target width.
$mdw = 6
$targetPixels = [Math]::Round($targetStoredWidth * $mdw)
$columnWidthInput = ($targetPixels - 5) / $mdw

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $columnWidthInput
$overview.Columns.Item(6).ColumnWidth = $columnWidthInput

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $columnWidthInput

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $columnWidthInput
